$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowValues {
    param($row, $values)
    $col = 3  # column C
    foreach ($v in $values) {
        $ws.Cells.Item($row, $col).Value = $v
        $col = $col + 1
    }
}

# Rows 5-7: convert "OOR" inline-string placeholders to computed numeric results (C:L)
Set-RowValues 5 @(5415.524044850007, 3281.95083062353, 1066.786607113238, 1066.786607113238, 56.76909753471365, 13332.73375124685, 983.5743079156682, 533.4728362364156, 196.196514182822, 19361.86197272529)
Set-RowValues 6 @(1446.220395474019, 687.3528187524578, 1066.786607113238, 1066.786607113238, 0.9336061902159818, 1186.664584156331, 983.5743079156682, 47.48113426565523, 31.12607746882464, 10862.32486602186)
Set-RowValues 7 @(1294.570534682059, 839.0026795444178, 1066.786607113238, 1066.786607113238, -1.199625539479206, 722.6160101329336, 983.5743079156682, 28.9135011339597, 29.16020922704595, 7060.511162914206)

# Rows 8-10: update already-numeric results with the new computed values
Set-RowValues 8 @(1215.30424107269, 918.2689731537868, 1066.786607113238, 1066.786607113238, -2.314650451073491, 480.0611516882645, 983.5743079156682, 19.20833258476407, 25.77879439013542, 5305.827915326061)
Set-RowValues 9 @(1172.594880786659, 960.9783334398184, 1066.786607113238, 1066.786607113238, -2.915435457016446, 349.3705092130081, 983.5743079156682, 13.97910435508354, 24.73314910841883, 4024.630623436303)
Set-RowValues 10 @(1151.309911825298, 982.2633024011792, 1066.786607113238, 1066.786607113238, -3.214847345143676, 284.2385041912438, 983.5743079156682, 11.37302550456456, 25.24742739549875, 3207.635248897907)

# Updated input parameters (column B) feeding the graphing calculations
$ws.Range("B14").Value = 1.53
$ws.Range("B15").Value = 0.2
$ws.Range("B16").Value = 0.012
$ws.Range("B17").Value = 0.5
$ws.Range("B18").Value = 95
$ws.Range("B20").Value = 2134
$ws.Range("B22").Value = 3.5
$ws.Range("B23").Value = 0.765
$ws.Range("B24").Value = 1.81
$ws.Range("B25").Value = 1.07
